# Auto-generated Excel COM-interop script to apply horarios update
$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")

# Cell updates (value changes / swaps)
$ws.Cells.Item(2, 1).Value = "Última actualización: 14:12:26"
$ws.Cells.Item(3, 1).Value = "Total filas: 204"
$ws.Cells.Item(56, 1).Value = "08:27:16"
$ws.Cells.Item(56, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(56, 4).Value = 50
$ws.Cells.Item(57, 1).Value = "07:38:39"
$ws.Cells.Item(57, 3).Value = "14_ABASTO"
$ws.Cells.Item(57, 4).Value = 99
$ws.Cells.Item(106, 1).Value = "10:05:51"
$ws.Cells.Item(106, 3).Value = "225_GOMEZ"
$ws.Cells.Item(106, 4).Value = 107
$ws.Cells.Item(108, 1).Value = "11:47:17"
$ws.Cells.Item(108, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(108, 4).Value = 5
$ws.Cells.Item(120, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(121, 3).Value = "15_ABASTO"
$ws.Cells.Item(133, 1).Value = "11:11:33"
$ws.Cells.Item(133, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(133, 4).Value = 84
$ws.Cells.Item(134, 1).Value = "11:34:59"
$ws.Cells.Item(134, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(134, 4).Value = 61
$ws.Cells.Item(146, 1).Value = "11:47:17"
$ws.Cells.Item(146, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(146, 4).Value = 76
$ws.Cells.Item(147, 1).Value = "11:34:59"
$ws.Cells.Item(147, 3).Value = "215C_EL PATO"
$ws.Cells.Item(147, 4).Value = 89
$ws.Cells.Item(158, 1).Value = "12:11:52"
$ws.Cells.Item(158, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(158, 4).Value = 74
$ws.Cells.Item(159, 1).Value = "11:47:17"
$ws.Cells.Item(159, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(159, 4).Value = 98
$ws.Cells.Item(197, 1).Value = "14:12:26"
$ws.Cells.Item(197, 2).Value = "15:17"
$ws.Cells.Item(197, 4).Value = 65
$ws.Cells.Item(198, 1).Value = "13:56:11"
$ws.Cells.Item(198, 2).Value = "15:18"
$ws.Cells.Item(198, 3).Value = "14_ABASTO"
$ws.Cells.Item(198, 4).Value = 82
$ws.Cells.Item(199, 1).Value = "14:12:26"
$ws.Cells.Item(199, 2).Value = "15:29"
$ws.Cells.Item(199, 3).Value = "10_OLMOS"
$ws.Cells.Item(199, 4).Value = 77
$ws.Cells.Item(200, 1).Value = "13:41:54"
$ws.Cells.Item(200, 2).Value = "15:32"
$ws.Cells.Item(200, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(200, 4).Value = 111
$ws.Cells.Item(201, 1).Value = "13:41:54"
$ws.Cells.Item(201, 2).Value = "15:33"
$ws.Cells.Item(201, 3).Value = "215C_EL PATO"
$ws.Cells.Item(201, 4).Value = 112
$ws.Cells.Item(202, 2).Value = "15:34"
$ws.Cells.Item(202, 3).Value = "215C_EL PATO"
$ws.Cells.Item(202, 4).Value = 98

# New rows appended
$ws.Cells.Item(203, 1).Value = "14:12:26"
$ws.Cells.Item(203, 2).Value = "15:36"
$ws.Cells.Item(203, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(203, 4).Value = 84
$ws.Cells.Item(203, 5).Value = "LP1912"
$ws.Cells.Item(204, 1).Value = "14:12:26"
$ws.Cells.Item(204, 2).Value = "15:41"
$ws.Cells.Item(204, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(204, 4).Value = 89
$ws.Cells.Item(204, 5).Value = "LP1912"
$ws.Cells.Item(205, 1).Value = "13:56:11"
$ws.Cells.Item(205, 2).Value = "15:42"
$ws.Cells.Item(205, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(205, 4).Value = 106
$ws.Cells.Item(205, 5).Value = "LP1912"
$ws.Cells.Item(206, 1).Value = "13:56:11"
$ws.Cells.Item(206, 2).Value = "15:53"
$ws.Cells.Item(206, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(206, 4).Value = 117
$ws.Cells.Item(206, 5).Value = "LP1912"
$ws.Cells.Item(207, 1).Value = "13:56:11"
$ws.Cells.Item(207, 2).Value = "15:53"
$ws.Cells.Item(207, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(207, 4).Value = 117
$ws.Cells.Item(207, 5).Value = "LP1912"
$ws.Cells.Item(208, 1).Value = "14:12:26"
$ws.Cells.Item(208, 2).Value = "15:56"
$ws.Cells.Item(208, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(208, 4).Value = 104
$ws.Cells.Item(208, 5).Value = "LP1912"
$ws.Cells.Item(209, 1).Value = "14:12:26"
$ws.Cells.Item(209, 2).Value = "16:05"
$ws.Cells.Item(209, 3).Value = "14_ABASTO"
$ws.Cells.Item(209, 4).Value = 113
$ws.Cells.Item(209, 5).Value = "LP1912"

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")

# Cell updates (value changes / swaps)
$ws.Cells.Item(2, 1).Value = "Última actualización: 14:12:26"

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")

# Cell updates (value changes / swaps)
$ws.Cells.Item(2, 1).Value = "Última actualización: 14:12:26"
$ws.Cells.Item(3, 1).Value = "Total filas: 29"
$ws.Cells.Item(19, 1).Value = "08:37:25"
$ws.Cells.Item(19, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(19, 4).Value = 113
$ws.Cells.Item(20, 1).Value = "08:52:50"
$ws.Cells.Item(20, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(20, 4).Value = 98

# New rows appended
$ws.Cells.Item(34, 1).Value = "14:12:26"
$ws.Cells.Item(34, 2).Value = "16:02"
$ws.Cells.Item(34, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(34, 4).Value = 110
$ws.Cells.Item(34, 5).Value = "L6203"
